# Commit: "Now has edit product from manage once done open manage and press
# stop edit" -- removing a product entry (row 14: "juice" / "22") from the
# product database sheet, e.g. after deleting/editing it from the "manage
# products" screen and pressing "stop edit". The rows below it (16, 17) keep
# their original row numbers, so this is a content clear of the row, not a
# row deletion/shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14:B14").ClearContents()
